# Update the division-practice answer table.
# The document contains a single 5-column table; only rows 1, 5, 9, 13, 17
# (1-based) hold data, each with 5 cells of "a÷b=c, d" text.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("19÷8=2, 3", "45÷8=5, 5", "85÷6=14, 1", "36÷5=7, 1", "66÷9=7, 3")
    5  = @("87÷8=10, 7", "22÷2=11, 0", "27÷8=3, 3", "28÷2=14, 0", "21÷8=2, 5")
    9  = @("59÷2=29, 1", "34÷5=6, 4", "23÷7=3, 2", "22÷7=3, 1", "43÷9=4, 7")
    13 = @("99÷5=19, 4", "55÷6=9, 1", "30÷9=3, 3", "25÷9=2, 7", "38÷9=4, 2")
    17 = @("35÷3=11, 2", "67÷9=7, 4", "64÷3=21, 1", "66÷4=16, 2", "98÷5=19, 3")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
